# Auto-generated Excel COM-interop edit script
# Applies the "New crime data collected" update to the CompStat weekly report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report title / date range)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value2 = "Volume 32   Number  5"
$ws.Range("C9").Value2 = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# ---------------------------------------------------------------------------
# Column E width adjustment (bestFit recalculated by the data update)
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 6.71

# ---------------------------------------------------------------------------
# Reference cells used to clone number-format/style when a cell's data type
# changes between "numeric" and "text placeholder" (e.g. "0" / "***.*").
# These reference cells are never themselves modified by this script.
#   style 13 -> General/text placeholder style  (sample: A14)
#   style 14 -> integer "#,##0" style            (sample: G15)
#   style 15 -> decimal "#,##0.0" (%chg) style    (sample: N14)
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Crime-statistics table updates (rows 14-30)
# ---------------------------------------------------------------------------
$c = $ws.Range("F14")
$c.NumberFormat = "@"
$c.Value2 = "0"
$ws.Range("A14").Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("M14")
$c.Value2 = 0
$ws.Range("N14").Copy()
$c.PasteSpecial(-4122)
$c.Value2 = 0
$ws.Range("C15").Value2 = 2
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = "0"
$ws.Range("A14").Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value2 = "***.*"
$ws.Range("A14").Copy()
$c.PasteSpecial(-4122)
$ws.Range("F15").Value2 = 4
$ws.Range("H15").Value2 = 100
$ws.Range("I15").Value2 = 6
$ws.Range("K15").Value2 = 200
$ws.Range("L15").Value2 = 500
$ws.Range("M15").Value2 = 200
$ws.Range("N15").Value2 = -25
$ws.Range("C16").Value2 = 6
$ws.Range("D16").Value2 = 5
$ws.Range("E16").Value2 = 20
$ws.Range("F16").Value2 = 17
$ws.Range("G16").Value2 = 19
$ws.Range("H16").Value2 = -10.526315789473
$ws.Range("I16").Value2 = 19
$ws.Range("J16").Value2 = 22
$ws.Range("K16").Value2 = -13.636363636363
$ws.Range("L16").Value2 = 18.75
$ws.Range("M16").Value2 = -17.391304347826
$ws.Range("N16").Value2 = -81.553398058252
$ws.Range("D17").Value2 = 9
$ws.Range("E17").Value2 = -55.555555555555
$ws.Range("F17").Value2 = 23
$ws.Range("G17").Value2 = 35
$ws.Range("H17").Value2 = -34.285714285714
$ws.Range("I17").Value2 = 32
$ws.Range("J17").Value2 = 43
$ws.Range("K17").Value2 = -25.581395348837
$ws.Range("L17").Value2 = -17.948717948717
$ws.Range("M17").Value2 = 6.666666666666
$ws.Range("N17").Value2 = -52.238805970149
$ws.Range("D18").Value2 = 5
$ws.Range("E18").Value2 = -40
$ws.Range("F18").Value2 = 8
$ws.Range("H18").Value2 = -20
$ws.Range("I18").Value2 = 11
$ws.Range("J18").Value2 = 12
$ws.Range("K18").Value2 = -8.333333333333
$ws.Range("L18").Value2 = -50
$ws.Range("M18").Value2 = -26.666666666666
$ws.Range("N18").Value2 = -87.209302325581
$ws.Range("C19").Value2 = 6
$ws.Range("D19").Value2 = 5
$ws.Range("E19").Value2 = 20
$ws.Range("G19").Value2 = 22
$ws.Range("H19").Value2 = 40.909090909090
$ws.Range("I19").Value2 = 36
$ws.Range("J19").Value2 = 26
$ws.Range("K19").Value2 = 38.461538461538
$ws.Range("L19").Value2 = 80
$ws.Range("M19").Value2 = 56.521739130434
$ws.Range("N19").Value2 = -28
$ws.Range("C20").Value2 = 2
$ws.Range("D20").Value2 = 3
$ws.Range("E20").Value2 = -33.333333333333
$ws.Range("F20").Value2 = 6
$ws.Range("G20").Value2 = 12
$ws.Range("H20").Value2 = -50
$ws.Range("I20").Value2 = 9
$ws.Range("J20").Value2 = 12
$ws.Range("K20").Value2 = -25
$ws.Range("L20").Value2 = 80
$ws.Range("M20").Value2 = 80
$ws.Range("N20").Value2 = -66.666666666666
$ws.Range("C21").Value2 = 23
$ws.Range("D21").Value2 = 27
$ws.Range("E21").Value2 = -14.814814814814
$ws.Range("F21").Value2 = 89
$ws.Range("G21").Value2 = 100
$ws.Range("H21").Value2 = -11
$ws.Range("I21").Value2 = 114
$ws.Range("J21").Value2 = 117
$ws.Range("K21").Value2 = -2.564102564102
$ws.Range("L21").Value2 = 10.679611650485
$ws.Range("M21").Value2 = 15.151515151515
$ws.Range("N21").Value2 = -67.052023121387
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value2 = "0"
$ws.Range("A14").Copy()
$c.PasteSpecial(-4122)
$ws.Range("F22").Value2 = 2
$ws.Range("H22").Value2 = 100
$c = $ws.Range("L22")
$c.Value2 = 200
$ws.Range("N14").Copy()
$c.PasteSpecial(-4122)
$c.Value2 = 200
$ws.Range("C23").Value2 = 3
$ws.Range("E23").Value2 = -40
$ws.Range("F23").Value2 = 15
$ws.Range("G23").Value2 = 20
$ws.Range("H23").Value2 = -25
$ws.Range("I23").Value2 = 19
$ws.Range("J23").Value2 = 25
$ws.Range("K23").Value2 = -24
$ws.Range("L23").Value2 = 0
$ws.Range("M23").Value2 = 72.727272727272
$ws.Range("C24").Value2 = 14
$ws.Range("D24").Value2 = 17
$ws.Range("E24").Value2 = -17.647058823529
$ws.Range("F24").Value2 = 71
$ws.Range("G24").Value2 = 62
$ws.Range("H24").Value2 = 14.516129032258
$ws.Range("I24").Value2 = 79
$ws.Range("J24").Value2 = 74
$ws.Range("K24").Value2 = 6.756756756756
$ws.Range("L24").Value2 = 2.597402597402
$ws.Range("M24").Value2 = 33.898305084745
$ws.Range("C25").Value2 = 6
$ws.Range("D25").Value2 = 3
$ws.Range("E25").Value2 = 100
$ws.Range("F25").Value2 = 13
$ws.Range("H25").Value2 = -35
$ws.Range("I25").Value2 = 14
$ws.Range("J25").Value2 = 20
$ws.Range("K25").Value2 = -30
$ws.Range("L25").Value2 = -26.315789473684
$ws.Range("C26").Value2 = 11
$ws.Range("D26").Value2 = 10
$ws.Range("E26").Value2 = 10
$ws.Range("F26").Value2 = 42
$ws.Range("G26").Value2 = 34
$ws.Range("H26").Value2 = 23.529411764705
$ws.Range("I26").Value2 = 47
$ws.Range("J26").Value2 = 46
$ws.Range("K26").Value2 = 2.173913043478
$ws.Range("L26").Value2 = -17.543859649122
$ws.Range("M26").Value2 = -38.157894736842
$ws.Range("C27").Value2 = 2
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = "0"
$ws.Range("A14").Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value2 = "***.*"
$ws.Range("A14").Copy()
$c.PasteSpecial(-4122)
$ws.Range("F27").Value2 = 4
$ws.Range("H27").Value2 = 33.333333333333
$ws.Range("I27").Value2 = 6
$ws.Range("K27").Value2 = 100
$ws.Range("L27").Value2 = 200
$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value2 = "0"
$ws.Range("A14").Copy()
$c.PasteSpecial(-4122)
$ws.Range("E28").Value2 = -100
$ws.Range("F28").Value2 = 3
$ws.Range("H28").Value2 = -25
$ws.Range("J28").Value2 = 5
$ws.Range("K28").Value2 = -20
$ws.Range("L28").Value2 = 33.333333333333
$c = $ws.Range("C29")
$c.Value2 = 1
$ws.Range("G15").Copy()
$c.PasteSpecial(-4122)
$c.Value2 = 1
$c = $ws.Range("F29")
$c.Value2 = 1
$ws.Range("G15").Copy()
$c.PasteSpecial(-4122)
$c.Value2 = 1
$ws.Range("H29").Value2 = -66.666666666666
$c = $ws.Range("I29")
$c.Value2 = 1
$ws.Range("G15").Copy()
$c.PasteSpecial(-4122)
$c.Value2 = 1
$ws.Range("K29").Value2 = -66.666666666666
$ws.Range("L29").Value2 = -66.666666666666
$ws.Range("M29").Value2 = -66.666666666666
$ws.Range("N29").Value2 = -92.857142857142
$c = $ws.Range("C30")
$c.Value2 = 1
$ws.Range("G15").Copy()
$c.PasteSpecial(-4122)
$c.Value2 = 1
$c = $ws.Range("F30")
$c.Value2 = 1
$ws.Range("G15").Copy()
$c.PasteSpecial(-4122)
$c.Value2 = 1
$ws.Range("H30").Value2 = -50
$c = $ws.Range("I30")
$c.Value2 = 1
$ws.Range("G15").Copy()
$c.PasteSpecial(-4122)
$c.Value2 = 1
$ws.Range("K30").Value2 = -50
$ws.Range("L30").Value2 = -50
$ws.Range("M30").Value2 = -50
$ws.Range("N30").Value2 = -90.909090909090
$excel.CutCopyMode = 0
